# Rename "Sheet1" to "Shrinkage" (matches new sheet tab/sheet.xml name).
$wb = $excel.ActiveWorkbook
$shrinkage = $wb.Worksheets.Item("Sheet1")
$shrinkage.Name = "Shrinkage"

# The GSM sheet was previously the active tab (scrolled all the way down,
# whole-sheet selected). Reset its view back to a normal state before
# handing focus over, so it no longer carries the stale tabSelected /
# topLeftCell / full-sheet selection state.
$gsm = $wb.Worksheets.Item("GSM")
$gsm.Activate()
$gsm.Range("A1").Select()

# Make the newly-renamed "Shrinkage" sheet the active/selected tab.
$shrinkage.Activate()
